$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN" (sheet1): add rows 27-38, dimension A1:G26 -> A1:G38
# ---------------------------------------------------------------------------
$amsin = $wb.Worksheets.Item("AMSIN")

$amsinRows = @(
    @{ r=27; a="2021-11-17"; b=44517.63177959491;  c="153rsgn";         d=89; e=86; f=3;  g=4.56 },
    @{ r=28; a="2021-11-18"; b=44518.44433387731;  c="153fnlrgsn";      d=89; e=89; f=0;  g=2.41 },
    @{ r=29; a="2021-12-03"; b=44533.72964;        c="154_scndrgsn";    d=89; e=89; f=0;  g=2.71 },
    @{ r=30; a="2021-12-06"; b=44536.45990143518;  c="154_fnlrgsn";     d=89; e=89; f=0;  g=2.8  },
    @{ r=31; a="2021-12-23"; b=44553.46343936343;  c="155_fnlrgsn";     d=89; e=89; f=0;  g=3.66 },
    @{ r=32; a="2022-01-03"; b=44564.64808984954;  c="lodash";          d=89; e=89; f=0;  g=3.04 },
    @{ r=33; a="2022-01-17"; b=44578.53799537037;  c="frstrgsn156";     d=89; e=62; f=27; g=8.18 },
    @{ r=34; a="2022-01-19"; b=44580.77628545139;  c="165_secondcyc";   d=89; e=89; f=0;  g=2.46 },
    @{ r=35; a="2022-01-20"; b=44581.4165858912;   c="156_fnlrsgn";     d=89; e=87; f=2;  g=3.22 },
    @{ r=36; a="2022-01-28"; b=44589.60213442129;  c="156audit";        d=89; e=89; f=0;  g=2.29 },
    @{ r=37; a="2022-02-07"; b=44599.79848405092;  c="secondcycle_157"; d=89; e=89; f=0;  g=2.71 },
    @{ r=38; a="2022-02-08"; b=44600.43772320602;  c="157_fnl";         d=89; e=89; f=0;  g=2.74 }
)

# One-time: register the custom number format "yyyy-mm-dd h:mm:ss" (numFmtId 166)
# in the style table, matching the upstream edit's styles.xml, by applying it to
# a scratch cell and then restoring its real formatting from a template cell.
$amsin.Range("B27").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Column A holds literal "yyyy-mm-dd" strings (not real dates) -- force text so
# Excel doesn't silently convert them to date serials.
$amsinTextRange = $amsin.Range("A27:A38")
$amsinTextRange.NumberFormat = "@"
foreach ($row in $amsinRows) {
    $amsin.Cells.Item($row.r, 1).Value = $row.a
}
# Re-apply the normal (non-text) cell format from an existing templated row so
# the style index matches the rest of column A (style 7), then restore values.
$amsin.Range("A26").Copy()
$amsinTextRange.PasteSpecial(-4122) | Out-Null

foreach ($row in $amsinRows) {
    $r = $row.r
    # Column B: numeric timestamp, copy format from an existing "B" cell that
    # already carries the correct datetime style (style 14), then set the value.
    $amsin.Range("B26").Copy($amsin.Range("B$r"))
    $amsin.Cells.Item($r, 2).Value = $row.b

    $amsin.Cells.Item($r, 3).Value = $row.c
    $amsin.Cells.Item($r, 4).Value = $row.d
    $amsin.Cells.Item($r, 5).Value = $row.e
    $amsin.Cells.Item($r, 6).Value = $row.f
    $amsin.Cells.Item($r, 7).Value = $row.g
}

# ---------------------------------------------------------------------------
# Sheet "BETA" (sheet2): fix row 13, fill rows 14-16, add rows 17-18
# dimension A1:G16 -> A1:G18
# ---------------------------------------------------------------------------
$beta = $wb.Worksheets.Item("BETA")

# Row 13: only the B13 timestamp precision + formatting of the other cells changes.
$beta.Range("B12").Copy($beta.Range("A13"))
$beta.Cells.Item(13, 1).Value = "2021-10-28"
$beta.Range("B12").Copy($beta.Range("C13"))
$beta.Cells.Item(13, 3).Value = "152_betachgs"
$beta.Range("B12").Copy($beta.Range("E13"))
$beta.Cells.Item(13, 5).Value = 89
$beta.Range("B12").Copy($beta.Range("F13"))
$beta.Cells.Item(13, 6).Value = 0
$beta.Range("B12").Copy($beta.Range("G13"))
$beta.Cells.Item(13, 7).Value = 4.23
$beta.Cells.Item(13, 2).Value = 44497.70566909722

$betaRows1416 = @(
    @{ r=14; a="2021-11-18"; b=44518.5958771875;   c="153_beta"; d=89; e=89; f=0; g=2.94 },
    @{ r=15; a="2021-12-06"; b=44536.56433603009;  c="154_beta"; d=89; e=87; f=2; g=2.87 },
    @{ r=16; a="2021-12-23"; b=44553.53586863426;  c="155_beta"; d=89; e=89; f=0; g=2.57 }
)

foreach ($row in $betaRows1416) {
    $r = $row.r
    $beta.Range("A12").Copy($beta.Range("A$r"))
    $beta.Range("C12").Copy($beta.Range("C$r"))
    $beta.Range("D12").Copy($beta.Range("D$r"))
    $beta.Range("E12").Copy($beta.Range("E$r"))
    $beta.Range("F12").Copy($beta.Range("F$r"))
    $beta.Range("G12").Copy($beta.Range("G$r"))
    $beta.Range("B12").Copy($beta.Range("B$r"))

    $beta.Cells.Item($r, 1).Value = $row.a
    $beta.Cells.Item($r, 2).Value = $row.b
    $beta.Cells.Item($r, 3).Value = $row.c
    $beta.Cells.Item($r, 4).Value = $row.d
    $beta.Cells.Item($r, 5).Value = $row.e
    $beta.Cells.Item($r, 6).Value = $row.f
    $beta.Cells.Item($r, 7).Value = $row.g
}

# Row 16's "Total Cases" column uses the plain style (7), not the highlighted
# one (12) used by D11/D12 -- re-base it off a plain-style column-D cell.
$beta.Range("E16").Copy($beta.Range("D16"))
$beta.Cells.Item(16, 4).Value = 89

# Row 17 (brand new row, style 14 on the timestamp column).
$beta.Range("A12").Copy($beta.Range("A17"))
$beta.Cells.Item(17, 1).Value = "2022-01-20"
$beta.Range("C12").Copy($beta.Range("C17"))
$beta.Cells.Item(17, 3).Value = "156_beta"
$beta.Range("E12").Copy($beta.Range("E17"))
$beta.Cells.Item(17, 5).Value = 89
$beta.Range("E12").Copy($beta.Range("F17"))
$beta.Cells.Item(17, 6).Value = 0
$beta.Range("E12").Copy($beta.Range("G17"))
$beta.Cells.Item(17, 7).Value = 2.46
$amsin.Range("B26").Copy($beta.Range("B17"))
$beta.Cells.Item(17, 2).Value = 44581.55118877315

# Row 18 (brand new row, default style throughout except the style-14 timestamp).
$beta.Cells.Item(18, 1).Value = "2022-02-08"
$beta.Cells.Item(18, 3).Value = "157_beta"
$beta.Cells.Item(18, 4).Value = 89
$beta.Cells.Item(18, 5).Value = 89
$beta.Cells.Item(18, 6).Value = 0
$beta.Cells.Item(18, 7).Value = 3.05
$amsin.Range("B26").Copy($beta.Range("B18"))
$beta.Cells.Item(18, 2).Value = 44600.66145045477

# ---------------------------------------------------------------------------
# Sheet "AMS" (sheet3): add rows 19-23, dimension A1:G18 -> A1:G23
# ---------------------------------------------------------------------------
$ams = $wb.Worksheets.Item("AMS")

$amsRows = @(
    @{ r=19; a="2021-11-23"; b=44523.45710248843; c="153_live";  d=89; e=88; f=1; g=2.9  },
    @{ r=20; a="2021-12-06"; b=44536.88133276621; c="154_live";  d=89; e=87; f=2; g=3.25 },
    @{ r=21; a="2021-12-23"; b=44553.81873172454; c="155_live";  d=89; e=86; f=3; g=4.31 },
    @{ r=22; a="2021-12-27"; b=44557.71194859954; c="155hftfxx"; d=89; e=89; f=0; g=2.35 },
    @{ r=23; a="2022-01-20"; b=44581.83430385416; c="156_live";  d=89; e=89; f=0; g=3.3  }
)

foreach ($row in $amsRows) {
    $r = $row.r
    $ams.Range("A18").Copy($ams.Range("A$r"))
    $ams.Range("C18").Copy($ams.Range("C$r"))
    $ams.Range("D18").Copy($ams.Range("D$r"))
    $ams.Range("E18").Copy($ams.Range("E$r"))
    $ams.Range("F18").Copy($ams.Range("F$r"))
    $ams.Range("G18").Copy($ams.Range("G$r"))
    $ams.Range("B18").Copy($ams.Range("B$r"))

    $ams.Cells.Item($r, 1).Value = $row.a
    $ams.Cells.Item($r, 2).Value = $row.b
    $ams.Cells.Item($r, 3).Value = $row.c
    $ams.Cells.Item($r, 4).Value = $row.d
    $ams.Cells.Item($r, 5).Value = $row.e
    $ams.Cells.Item($r, 6).Value = $row.f
    $ams.Cells.Item($r, 7).Value = $row.g
}

Write-Output "done"
